$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "costs" worksheet after "production" (last sheet),
#    giving it sheetId 7 / rId7 at the end of the tab strip.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "costs"

# ------------------------------------------------------------------
# 2. Populate the template content. Writing the new shared strings in
#    this order (CostIncome, Amount, Rule, Project X) matches the
#    order they were first introduced in the source workbook.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Template"
$ws.Range("B1").Value = "CostIncome"

$ws.Range("A2").Value = "Date"
$ws.Range("C2").Value = "Amount"
$ws.Range("B2").Value = "Rule"

$ws.Range("A3").Value = 46387
$ws.Range("A3").NumberFormat = "m/d/yy"
$ws.Range("B3").Value = "Project X"
$ws.Range("C3").Value = -10000

# ------------------------------------------------------------------
# 3. Column widths for the new sheet.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.67
$ws.Columns.Item(2).ColumnWidth = 16.5

# ------------------------------------------------------------------
# 4. Update selections on the other sheets that moved (interest
#    rates, tax, audit, production) - done after the sheet-add so the
#    add's own view-normalisation pass doesn't clobber them again.
# ------------------------------------------------------------------
$wsInterest = $wb.Worksheets.Item("interest rates")
$wsInterest.Range("C38").Select()

$wsTax = $wb.Worksheets.Item("tax")
$wsTax.Range("D25").Select()

$wsAudit = $wb.Worksheets.Item("audit")
$wsAudit.Range("D24").Select()

$wsProduction = $wb.Worksheets.Item("production")
$wsProduction.Range("B14").Select()

# ------------------------------------------------------------------
# 5. Make "costs" the active sheet/tab with C3 selected - it becomes
#    the last-saved tab (activeTab) and carries tabSelected="1".
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("C3").Select()
